$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price-report row was inserted at row 169 (Femacal de La
# Calera - Albahaca), pushing the existing rows 169-246 down to 170-247.
$ws.Rows.Item(169).Insert()

$ws.Cells.Item(169, 1).Value = 3
$ws.Cells.Item(169, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(169, 3).Value = "Coquimbo"
$ws.Cells.Item(169, 4).Value = 44992
$ws.Cells.Item(169, 5).Value = 5
$ws.Cells.Item(169, 6).Value = 100112052
$ws.Cells.Item(169, 7).Value = "Albahaca"
$ws.Cells.Item(169, 8).Value = "Sin especificar"
$ws.Cells.Item(169, 9).Value = "Primera"
$ws.Cells.Item(169, 10).Value = 78
$ws.Cells.Item(169, 11).Value = 5000
$ws.Cells.Item(169, 12).Value = 5000
$ws.Cells.Item(169, 13).Value = 5000
$ws.Cells.Item(169, 14).Value = "`$/docena de matas"
$ws.Cells.Item(169, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(169, 16).Value = 833
$ws.Cells.Item(169, 17).Value = 6
$ws.Cells.Item(169, 18).Value = "Hortaliza"
